# Update column G ("K" - strikeouts) values on Sheet1 to reflect the
# regenerated save_data (commit: "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals").
#
# Only column G (rows 2-17 and 19) changes; row 18 stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 2
    14 = 5
    15 = 0
    16 = 3
    17 = 2
    19 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
